$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.37000000000053
$ws.Range("H2").Value = [double]"1.831295710721908e-16"
$ws.Range("K2").Value = 49.60593769851756
$ws.Range("L2").Value = "[42.8128269249299, 56.399048472105214]"
$ws.Range("O2").Value = 1.364816027685656
$ws.Range("P2").Value = "[1.2138686329185786, 1.515763422452733]"
$ws.Range("S2").Value = 55.53725304696597
$ws.Range("T2").Value = "[51.378061343722685, 59.696444750209245]"
$ws.Range("W2").Value = 19.85919919919961
$ws.Range("X2").Value = 19.24970970971011
$ws.Range("Y2").Value = 20.46868868868911

# Row 3 updates
$ws.Range("E3").Value = 25.49000000000055
$ws.Range("H3").Value = [double]"1.831295710721908e-16"
# I3 goes from a numeric p_reject value to a blank/empty text cell (like I2).
# A plain Value = "" clears the cell outright, so prime it as text via the
# leading apostrophe (Excel's literal-text marker) then restore the
# original (unmarked) style so no stray formatting is left behind.
$i3Style = $ws.Range("I3").Style
$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = $i3Style
$ws.Range("K3").Value = 55.32752231041231
$ws.Range("L3").Value = "[48.33881316001128, 62.31623146081335]"
$ws.Range("O3").Value = 1.352237078121733
$ws.Range("P3").Value = "[1.2138686329185786, 1.490605523324887]"
$ws.Range("S3").Value = 57.76601424199836
$ws.Range("T3").Value = "[53.63233692098918, 61.89969156300754]"
$ws.Range("W3").Value = 20.00416416416459
$ws.Range("X3").Value = 19.44282282282324
$ws.Range("Y3").Value = 20.56550550550595
